$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.905.14'
$ws.Range("E2").Value = '  -1.70%  '

$ws.Range("D3").Value = '1.833.51'
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.00'
$ws.Range("E5").Value = '  +0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6886'
$ws.Range("E6").Value = '  -2.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07686'
$ws.Range("E8").Value = '  -2.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3053'
$ws.Range("E9").Value = '  -2.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.39'
$ws.Range("E10").Value = '  -4.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07801'
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").Value = '1.827.10'
$ws.Range("E12").Value = '  -3.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.090'
$ws.Range("E13").Value = '  -1.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.39'
$ws.Range("E14").Value = '  -3.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6818'
$ws.Range("E15").Value = '  -2.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.446'
$ws.Range("E16").Value = '  -1.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008296'
$ws.Range("E17").Value = '  -1.38%  '

$ws.Range("D18").Value = '28.895.89'
$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.60'
$ws.Range("E19").Value = '  -3.81%  '

$ws.Range("D20").Value = '2.074.44'
$ws.Range("E20").Value = '  -3.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.73'
$ws.Range("E21").Value = '  -2.95%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.471'
$ws.Range("E23").Value = '  -2.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9998'
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1476'
$ws.Range("E25").Value = '  -4.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.50'
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("E27").Value = '  -2.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.19'
$ws.Range("E28").Value = '  -3.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.544'
$ws.Range("E29").Value = '  +2.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.213'
$ws.Range("E30").Value = '  -2.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.155'
$ws.Range("E31").Value = '  -2.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.183'
$ws.Range("E32").Value = '  -2.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05114'
$ws.Range("E33").Value = '  -2.92%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7659'
$ws.Range("E34").Value = '  +1.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.845'
$ws.Range("E35").Value = '  -2.86%  '

$ws.Range("E36").Value = '  -3.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.690'
$ws.Range("E37").Value = '  -0.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01850'
$ws.Range("E38").Value = '  -1.47%  '

$ws.Range("D39").Value = '1.217.98'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.701'
$ws.Range("E40").Value = '  -2.57%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9406'
$ws.Range("E41").Value = '  +4.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.64'
$ws.Range("E42").Value = '  -0.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9995'
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.683'
$ws.Range("E44").Value = '  -5.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("E45").Value = '  -4.35%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.531'
$ws.Range("E46").Value = '  -0.75%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5167'
$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("D48").Value = '1.975.16'
$ws.Range("E48").Value = '  -3.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.13'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.749'
$ws.Range("E50").Value = '  -3.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4185'
$ws.Range("E51").Value = '  -2.66%  '
